$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 463565.66
$ws.Range("I28").Value = 694704.75
$ws.Range("J28").Value = 1287.5
$ws.Range("K28").Value = 694704.75
$ws.Range("L28").Value = 1287.5
$ws.Range("M28").Value = -694219.75
$ws.Range("N28").Value = -2257.5

$ws.Range("H112").Value = 17046480
$ws.Range("J112").Value = 19481664
$ws.Range("L112").Value = 58444992
$ws.Range("N112").Value = -58447208

$ws.Range("H129").Value = 1008.96295
$ws.Range("I129").Value = 441.66666
$ws.Range("J129").Value = 1292.6111
$ws.Range("K129").Value = 1324.99998
$ws.Range("L129").Value = 3877.8333
$ws.Range("M129").Value = 3675.00002
$ws.Range("N129").Value = -13877.8333

$ws.Range("H135").Value = 1220.875
$ws.Range("I135").Value = 1312.25
$ws.Range("J135").Value = 764
$ws.Range("K135").Value = 11810.25
$ws.Range("L135").Value = 6876
$ws.Range("M135").Value = -9275.25
$ws.Range("N135").Value = -11946

$ws.Range("H137").Value = 25642320
$ws.Range("I137").Value = 33334214
$ws.Range("J137").Value = 2675.889
$ws.Range("K137").Value = 100002642
$ws.Range("L137").Value = 8027.667
$ws.Range("M137").Value = -100000092
$ws.Range("N137").Value = -13127.667

$ws.Range("H138").Value = 2005.2828
$ws.Range("I138").Value = 749.6
$ws.Range("J138").Value = 2323.1772
$ws.Range("K138").Value = 2248.8
$ws.Range("L138").Value = 6969.5316
$ws.Range("M138").Value = 2891.2
$ws.Range("N138").Value = -17249.5316

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2388.96
$ws.Range("I61").Value = 1209.2
$ws.Range("K61").Value = 1209.2
$ws.Range("M61").Value = -997.2

$ws.Range("H74").Value = 4397.8716
$ws.Range("I74").Value = 1180.56
$ws.Range("J74").Value = 10143.071
$ws.Range("K74").Value = 1180.56
$ws.Range("L74").Value = 10143.071
$ws.Range("M74").Value = -306.5599999999999
$ws.Range("N74").Value = -11891.071

$ws.Range("H77").Value = 4397.8716
$ws.Range("I77").Value = 1180.56
$ws.Range("J77").Value = 10143.071
$ws.Range("K77").Value = 5902.799999999999
$ws.Range("L77").Value = 50715.355
$ws.Range("M77").Value = -1534.799999999999
$ws.Range("N77").Value = -59451.355

$ws.Range("H122").Value = 1856.7333
$ws.Range("I122").Value = 1555.3
$ws.Range("J122").Value = 2459.6
$ws.Range("K122").Value = 4665.9
$ws.Range("L122").Value = 7378.799999999999
$ws.Range("M122").Value = -2215.9
$ws.Range("N122").Value = -12278.8

$ws.Range("H132").Value = 2943.3096
$ws.Range("I132").Value = 2525.9092
$ws.Range("J132").Value = 4473.778
$ws.Range("K132").Value = 7577.7276
$ws.Range("L132").Value = 13421.334
$ws.Range("M132").Value = -5047.7276
$ws.Range("N132").Value = -18481.334

$ws.Range("H136").Value = 2388.96
$ws.Range("I136").Value = 1209.2
$ws.Range("K136").Value = 3627.6
$ws.Range("M136").Value = -1077.6

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2706.025
$ws.Range("I134").Value = 2106.1667
$ws.Range("J134").Value = 4505.6
$ws.Range("K134").Value = 6318.500100000001
$ws.Range("L134").Value = 13516.8
$ws.Range("M134").Value = -3783.500100000001
$ws.Range("N134").Value = -18586.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1187.9722
$ws.Range("I31").Value = 969.75
$ws.Range("J31").Value = 1951.75
$ws.Range("K31").Value = 969.75
$ws.Range("L31").Value = 1951.75
$ws.Range("M31").Value = -674.75
$ws.Range("N31").Value = -2541.75

$ws.Range("H34").Value = 1187.9722
$ws.Range("I34").Value = 969.75
$ws.Range("J34").Value = 1951.75
$ws.Range("K34").Value = 969.75
$ws.Range("L34").Value = 1951.75
$ws.Range("M34").Value = -767.75
$ws.Range("N34").Value = -2355.75

$ws.Range("H58").Value = 2189.2
$ws.Range("I58").Value = 1566.4706
$ws.Range("K58").Value = 1566.4706
$ws.Range("M58").Value = -1363.4706

$ws.Range("H132").Value = 2158.0286
$ws.Range("I132").Value = 1789.2593
$ws.Range("K132").Value = 5367.7779
$ws.Range("M132").Value = -2837.7779

$ws.Range("H134").Value = 3068.6785
$ws.Range("I134").Value = 1403.5883
$ws.Range("J134").Value = 5642
$ws.Range("K134").Value = 4210.7649
$ws.Range("L134").Value = 16926
$ws.Range("M134").Value = -1675.7649
$ws.Range("N134").Value = -21996

$ws.Range("H136").Value = 2189.2
$ws.Range("I136").Value = 1566.4706
$ws.Range("K136").Value = 4699.4118
$ws.Range("M136").Value = -2149.4118

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 898.63635
$ws.Range("I5").Value = 503.8421
$ws.Range("J5").Value = 1434.4286
$ws.Range("K5").Value = 1511.5263
$ws.Range("L5").Value = 4303.2858
$ws.Range("M5").Value = -1399.5263
$ws.Range("N5").Value = -4527.2858

$ws.Range("H63").Value = 2544.4443
$ws.Range("J63").Value = 3000
$ws.Range("L63").Value = 9000
$ws.Range("N63").Value = -10498

$ws.Range("H64").Value = 2870
$ws.Range("I64").Value = 1480
$ws.Range("K64").Value = 4440
$ws.Range("M64").Value = -4170

$ws.Range("H66").Value = 2544.4443
$ws.Range("J66").Value = 3000
$ws.Range("L66").Value = 27000
$ws.Range("N66").Value = -34488

$ws.Range("H67").Value = 2870
$ws.Range("I67").Value = 1480
$ws.Range("K67").Value = 4440
$ws.Range("M67").Value = -3504

$ws.Range("H70").Value = 2679.5
$ws.Range("I70").Value = 1198.75
$ws.Range("J70").Value = 3666.6667
$ws.Range("K70").Value = 3596.25
$ws.Range("L70").Value = 11000.0001
$ws.Range("M70").Value = -3281.25
$ws.Range("N70").Value = -11630.0001

$ws.Range("H73").Value = 2679.5
$ws.Range("I73").Value = 1198.75
$ws.Range("J73").Value = 3666.6667
$ws.Range("K73").Value = 3596.25
$ws.Range("L73").Value = 11000.0001
$ws.Range("M73").Value = -2504.25
$ws.Range("N73").Value = -13184.0001

$ws.Range("H87").Value = 13230.117
$ws.Range("I87").Value = 4992
$ws.Range("K87").Value = 14976
$ws.Range("M87").Value = -13728

$ws.Range("H90").Value = 13230.117
$ws.Range("I90").Value = 4992
$ws.Range("K90").Value = 44928
$ws.Range("M90").Value = -38688

$ws.Range("H114").Value = 1435.5
$ws.Range("J114").Value = 1680.091
$ws.Range("L114").Value = 5040.272999999999
$ws.Range("N114").Value = -11548.273

$ws.Range("H129").Value = 1486.52
$ws.Range("J129").Value = 1789.5385
$ws.Range("L129").Value = 5368.6155
$ws.Range("N129").Value = -15368.6155

$ws.Range("H131").Value = 1471.3768
$ws.Range("J131").Value = 1583.3064
$ws.Range("L131").Value = 4749.9192
$ws.Range("N131").Value = -14829.9192

$ws.Range("H135").Value = 898.63635
$ws.Range("I135").Value = 503.8421
$ws.Range("J135").Value = 1434.4286
$ws.Range("K135").Value = 4534.5789
$ws.Range("L135").Value = 12909.8574
$ws.Range("M135").Value = -1999.5789
$ws.Range("N135").Value = -17979.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 29010.8
$ws.Range("J93").Value = 29010.8
$ws.Range("L93").Value = 29010.8
$ws.Range("N93").Value = -32754.8

$ws.Range("H132").Value = 2478.0728
$ws.Range("I132").Value = 2185.8857
$ws.Range("J132").Value = 2989.4
$ws.Range("K132").Value = 6557.657099999999
$ws.Range("L132").Value = 8968.2
$ws.Range("M132").Value = -4027.657099999999
$ws.Range("N132").Value = -14028.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 11637.272
$ws.Range("I61").Value = 12083.333
$ws.Range("J61").Value = 11102
$ws.Range("K61").Value = 12083.333
$ws.Range("L61").Value = 11102
$ws.Range("M61").Value = -11881.333
$ws.Range("N61").Value = -11506

$ws.Range("H100").Value = 2721.3044
$ws.Range("I100").Value = 2198.889
$ws.Range("K100").Value = 2198.889
$ws.Range("M100").Value = -1657.889

$ws.Range("H113").Value = 11637.272
$ws.Range("I113").Value = 12083.333
$ws.Range("J113").Value = 11102
$ws.Range("K113").Value = 12083.333
$ws.Range("L113").Value = 11102
$ws.Range("M113").Value = -9913.333
$ws.Range("N113").Value = -15442

$ws.Range("H132").Value = 3586.3635
$ws.Range("I132").Value = 2897.111
$ws.Range("J132").Value = 4892.316
$ws.Range("K132").Value = 8691.332999999999
$ws.Range("L132").Value = 14676.948
$ws.Range("M132").Value = -6161.332999999999
$ws.Range("N132").Value = -19736.948

$ws.Range("H136").Value = 3082.5232
$ws.Range("I136").Value = 2244.2327
$ws.Range("J136").Value = 4721
$ws.Range("K136").Value = 6732.6981
$ws.Range("L136").Value = 14163
$ws.Range("M136").Value = -4182.6981
$ws.Range("N136").Value = -19263

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 41674144
$ws.Range("I132").Value = 50008670
$ws.Range("K132").Value = 150026010
$ws.Range("M132").Value = -150023480

$ws.Range("H136").Value = 8573884
$ws.Range("I136").Value = 11529603
$ws.Range("J136").Value = 2299.5
$ws.Range("K136").Value = 34588809
$ws.Range("L136").Value = 6898.5
$ws.Range("M136").Value = -34586259
$ws.Range("N136").Value = -11998.5
